# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.544.27"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.618.73"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "1.849.27"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "1.620.35"
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "27.554.13"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E24").Value = "  +6.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.98%  "
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("D33").Value = "1.440.27"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.939"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.55%  "
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0168"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "68.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.59%  "
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("E45").Value = "  -3.26%  "
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("D47").Value = "1.759.42"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.100"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.13%  "
